$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header style (bold, bordered, centered) used by the other headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for new columns I (I0) and J (IF), rows 2-35
$data = @(
    @(5,7),
    @(9,9),
    @(7,7),
    @(1,1),
    @(8,8),
    @(5,6),
    @(6,7),
    @(8,8),
    @(1,1),
    @(8,9),
    @(6,6),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(7,8),
    @(8,9),
    @(4,6),
    @(8,9),
    @(8,9),
    @(7,8),
    @(8,9),
    @(9,9),
    @(7,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(3,3),
    @(7,7),
    @(7,7),
    @(7,7)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
